$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").Value = 2.3772960904422913
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.4158372435670477
$ws.Range("E2").ClearContents()

# Row 3 updates
$ws.Range("B3").Value = 2.2032100745536449
$ws.Range("C3").Value = -0.75226107008933984
$ws.Range("D3").Value = 2.7060271389174742
$ws.Range("E3").Value = -1.5154658291482421

# Keep selection consistent with the recorded change
$ws.Range("B1:E3").Select()
